$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (shifts existing rows 7-17 down to 8-18)
$ws.Rows(7).Insert()

# Fill in the new row 7 with the Danish border-fence entry.
# Columns are set left-to-right so new shared-string entries are appended
# in the same order the source workbook recorded them.
$ws.Range("A7").Value = "DNK"
$ws.Range("B7").Value = "DEU"
$ws.Range("C7").Value = 2019
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "fencing"
$ws.Range("F7").Value = "disease"
$ws.Range("G7").Value = "68"
$ws.Range("H7").Value = "checked"
$ws.Range("I7").Value = "NYT (2018); DW (2018)"
$ws.Range("J7").Value = "https://www.nytimes.com/2018/10/24/world/europe/pig-disease-denmark-swine-fever.html"
$ws.Range("K7").Value = "https://www.dw.com/en/denmark-to-build-controversial-german-border-fence/a-45078064"

# The row insert shifts cell content/styles down automatically, but it
# leaves the existing hyperlink anchored at its old address (J13), so
# re-create it on the new J14 cell. Hyperlinks.Add() resets the cell's
# text + style as a side effect, so restore both afterwards.
$ws.Range("J13").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J14"), "https://www.theguardian.com/world/2017/aug/24/russia-lithuania-border-fence-kaliningrad-estonia-eston-kohver;", "", "", "https://www.theguardian.com/world/2017/aug/24/russia-lithuania-border-fence-kaliningrad-estonia-eston-kohver; ")
$ws.Range("J14").Style = "Link"
$ws.Range("J14").Value = "https://www.theguardian.com/world/2017/aug/24/russia-lithuania-border-fence-kaliningrad-estonia-eston-kohver; https://www.dw.com/en/lithuanias-fence-on-kaliningrad-border/av-39731926"

# New hyperlink turning the Estonia/Russia source URL (now on row 9) into
# a clickable link.
$ws.Hyperlinks.Add($ws.Range("J9"), "https://www.rferl.org/a/estonia-fence-russia/27212586.html")
$ws.Range("J9").Style = "Link"

# Move the active selection to mirror the saved workbook state.
$ws.Range("L29").Select()
